# Updating an image in the PPT
# 1) Bump the auto date footer shown on the slide master + every slide
#    layout from 2024-01-27 to 2024-01-29.
# 2) Thicken and recolor the border of the "Ingredient" rounded-rectangle
#    callouts (slides 1, 3 and 4) from a thin teal line to a heavier
#    dark/theme-colored line.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# --- Date placeholder on the slide master ---
for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $sh = $m.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "2024-01-29"
    }
}

# --- Date placeholder on every slide layout ---
for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
    $lay = $m.CustomLayouts.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $sh = $lay.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "2024-01-29"
        }
    }
}

# --- "Ingredient" rounded-rectangle border on slides 1, 3 and 4 ---
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "Ingredient") {
                $sh.Line.Weight = 1
                $sh.Line.ForeColor.ObjectThemeColor = 1
            }
        }
    }
}
